$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 6
$ws.Range("H6").Value = 214.33333
$ws.Range("I6").Value = 185
$ws.Range("J6").Value = 317
$ws.Range("K6").Value = 555
$ws.Range("L6").Value = 951
$ws.Range("M6").Value = -443
$ws.Range("N6").Value = -1175

# ALC row 48
$ws.Range("H48").Value = 575
$ws.Range("I48").Value = 528.5714
$ws.Range("J48").Value = 900
$ws.Range("K48").Value = 1585.7142
$ws.Range("L48").Value = 2700
$ws.Range("M48").Value = -1293.7142
$ws.Range("N48").Value = -3284

# ALC row 51
$ws.Range("H51").Value = 4997.769
$ws.Range("I51").Value = 4975
$ws.Range("K51").Value = 4975
$ws.Range("M51").Value = -4491

# ALC row 56
$ws.Range("H56").Value = 575
$ws.Range("I56").Value = 528.5714
$ws.Range("J56").Value = 900
$ws.Range("K56").Value = 1585.7142
$ws.Range("L56").Value = 2700
$ws.Range("M56").Value = -1051.7142
$ws.Range("N56").Value = -3768

# ALC row 111
$ws.Range("H111").Value = 1794.6666
$ws.Range("I111").Value = 1979.375
$ws.Range("K111").Value = 5938.125
$ws.Range("M111").Value = -2871.125

# ALC row 113
$ws.Range("H113").Value = 8490.625
$ws.Range("J113").Value = 8154.1665
$ws.Range("L113").Value = 8154.1665
$ws.Range("N113").Value = -14662.1665

# ALC row 119
$ws.Range("H119").Value = 1824.8
$ws.Range("J119").Value = 1824.8
$ws.Range("L119").Value = 5474.4
$ws.Range("N119").Value = -15150.4

# ALC row 132
$ws.Range("H132").Value = 4952.037
$ws.Range("I132").Value = 4952.037
$ws.Range("K132").Value = 14856.111
$ws.Range("M132").Value = -12326.111

$ws = $wb.Worksheets.Item("ARM")
# ARM row 5
$ws.Range("H5").Value = 681.5
$ws.Range("I5").Value = 681.5
$ws.Range("K5").Value = 681.5
$ws.Range("M5").Value = -569.5

# ARM row 45
$ws.Range("H45").Value = 1250.8889
$ws.Range("I45").Value = 741.2
$ws.Range("K45").Value = 741.2
$ws.Range("M45").Value = -364.2

# ARM row 63
$ws.Range("H63").Value = 2138.4
$ws.Range("J63").Value = 2174.5
$ws.Range("L63").Value = 2174.5
$ws.Range("N63").Value = -3546.5

# ARM row 66
$ws.Range("H66").Value = 2138.4
$ws.Range("J66").Value = 2174.5
$ws.Range("L66").Value = 10872.5
$ws.Range("N66").Value = -17736.5

# ARM row 69
$ws.Range("H69").Value = 199729.5
$ws.Range("J69").Value = 199729.5
$ws.Range("L69").Value = 199729.5
$ws.Range("N69").Value = -201227.5

# ARM row 72
$ws.Range("H72").Value = 199729.5
$ws.Range("J72").Value = 199729.5
$ws.Range("L72").Value = 599188.5
$ws.Range("N72").Value = -606676.5

# ARM row 74
$ws.Range("H74").Value = 23810920
$ws.Range("I74").Value = 31250938
$ws.Range("K74").Value = 31250938
$ws.Range("M74").Value = -31250064

# ARM row 77
$ws.Range("H77").Value = 23810920
$ws.Range("I77").Value = 31250938
$ws.Range("K77").Value = 156254690
$ws.Range("M77").Value = -156250322

# ARM row 88
$ws.Range("H88").Value = 52503
$ws.Range("I88").Value = 100006
$ws.Range("J88").Value = 5000
$ws.Range("K88").Value = 100006
$ws.Range("L88").Value = 5000
$ws.Range("M88").Value = -99600
$ws.Range("N88").Value = -5812

# ARM row 91
$ws.Range("H91").Value = 52503
$ws.Range("I91").Value = 100006
$ws.Range("J91").Value = 5000
$ws.Range("K91").Value = 100006
$ws.Range("L91").Value = 5000
$ws.Range("M91").Value = -98602
$ws.Range("N91").Value = -7808

$ws = $wb.Worksheets.Item("BSM")
# BSM row 4
$ws.Range("H4").Value = 681.5
$ws.Range("I4").Value = 681.5
$ws.Range("K4").Value = 681.5
$ws.Range("M4").Value = -566.5

# BSM row 20
$ws.Range("H20").Value = 1804.32
$ws.Range("J20").Value = 3140.4
$ws.Range("L20").Value = 3140.4
$ws.Range("N20").Value = -3634.4

# BSM row 82
$ws.Range("H82").Value = 10623
$ws.Range("J82").Value = 46883
$ws.Range("L82").Value = 46883
$ws.Range("N82").Value = -47649

# BSM row 85
$ws.Range("H85").Value = 10623
$ws.Range("J85").Value = 46883
$ws.Range("L85").Value = 46883
$ws.Range("N85").Value = -49535

# BSM row 105
$ws.Range("H105").Value = 5420.7676
$ws.Range("I105").Value = 6439.8184
$ws.Range("J105").Value = 4353.1904
$ws.Range("K105").Value = 6439.8184
$ws.Range("L105").Value = 4353.1904
$ws.Range("M105").Value = -4692.8184
$ws.Range("N105").Value = -7847.1904

$ws = $wb.Worksheets.Item("CRP")
# CRP row 22
$ws.Range("H22").Value = 9051.385
$ws.Range("I22").Value = 11407.777
$ws.Range("K22").Value = 11407.777
$ws.Range("M22").Value = -11057.777

# CRP row 99
$ws.Range("H99").Value = 8295.799999999999
$ws.Range("I99").Value = 8106.4443
$ws.Range("K99").Value = 8106.4443
$ws.Range("M99").Value = -6608.4443

# CRP row 126
$ws.Range("H126").Value = 8295.799999999999
$ws.Range("I126").Value = 8106.4443
$ws.Range("K126").Value = 24319.3329
$ws.Range("M126").Value = -21849.3329

$ws = $wb.Worksheets.Item("CUL")
# CUL row 61
$ws.Range("H61").Value = 172.125
$ws.Range("I61").Value = 196.92308
$ws.Range("K61").Value = 590.76924
$ws.Range("M61").Value = -375.76924

# CUL row 80
$ws.Range("H80").Value = 854.4
$ws.Range("I80").Value = 499.5
$ws.Range("K80").Value = 1498.5
$ws.Range("M80").Value = -562.5

# CUL row 83
$ws.Range("H83").Value = 854.4
$ws.Range("I83").Value = 499.5
$ws.Range("K83").Value = 4495.5
$ws.Range("M83").Value = 184.5

# CUL row 98
$ws.Range("H98").Value = 531.3333
$ws.Range("J98").Value = 521.2222
$ws.Range("L98").Value = 1563.6666
$ws.Range("N98").Value = -4559.6666

# CUL row 107
$ws.Range("H107").Value = 1516
$ws.Range("J107").Value = 1516
$ws.Range("L107").Value = 4548
$ws.Range("N107").Value = -8388

$ws = $wb.Worksheets.Item("GSM")
# GSM row 102
$ws.Range("H102").Value = 2850
$ws.Range("I102").Value = 1200
$ws.Range("K102").Value = 1200
$ws.Range("M102").Value = 422

# GSM row 122
$ws.Range("H122").Value = 50002250
$ws.Range("J122").Value = 100002650
$ws.Range("L122").Value = 300007950
$ws.Range("N122").Value = -300012850

$ws = $wb.Worksheets.Item("LTW")
# LTW row 22
$ws.Range("H22").Value = 1234.4138
$ws.Range("I22").Value = 911.4375
$ws.Range("J22").Value = 1631.9231
$ws.Range("K22").Value = 911.4375
$ws.Range("L22").Value = 1631.9231
$ws.Range("M22").Value = -616.4375
$ws.Range("N22").Value = -2221.9231

# LTW row 27
$ws.Range("H27").Value = 1234.4138
$ws.Range("I27").Value = 911.4375
$ws.Range("J27").Value = 1631.9231
$ws.Range("K27").Value = 911.4375
$ws.Range("L27").Value = 1631.9231
$ws.Range("M27").Value = -804.4375
$ws.Range("N27").Value = -1845.9231

# LTW row 40
$ws.Range("H40").Value = 5130.207
$ws.Range("I40").Value = 4791.6
$ws.Range("K40").Value = 4791.6
$ws.Range("M40").Value = -4655.6

# LTW row 132
$ws.Range("H132").Value = 5497.032
$ws.Range("I132").Value = 2591.7144
$ws.Range("K132").Value = 7775.1432
$ws.Range("M132").Value = -5245.1432

$ws = $wb.Worksheets.Item("WVR")
# WVR row 68
$ws.Range("H68").Value = 51250.25
$ws.Range("J68").Value = 42500
$ws.Range("L68").Value = 42500
$ws.Range("N68").Value = -44122

# WVR row 71
$ws.Range("H71").Value = 51250.25
$ws.Range("J71").Value = 42500
$ws.Range("L71").Value = 127500
$ws.Range("N71").Value = -135612

# WVR row 132
$ws.Range("H132").Value = 1987.1177
$ws.Range("I132").Value = 1552.0667
$ws.Range("K132").Value = 4656.2001
$ws.Range("M132").Value = -2126.2001
